# Cadastro.xlsx - refactor commit:
#  - updates the sample record on the "Cadastro" sheet (sobrenome/email/nome)
#  - resets the selection on the "Lupa" sheet back to A2
#  - makes "Cadastro" the active/selected sheet again

$wb = $excel.ActiveWorkbook
$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsLupa     = $wb.Worksheets.Item("Lupa")

# Update the sample registration row (row 2) on the Cadastro sheet.
# Order matters for shared-string layout: B2/C2 are written before A2 so the
# new strings land in the same append order Excel itself produced.
$wsCadastro.Range("B2").Value = "Sales"
$wsCadastro.Range("C2").Value = "antsa@email.com"
$wsCadastro.Range("A2").Value = "Jefrey"

# Reset the Lupa sheet's selection to A2 (it previously pointed at A6).
$wsLupa.Activate()
$wsLupa.Range("A2").Select()

# Finish with Cadastro as the active/selected sheet.
$wsCadastro.Activate()
$wsCadastro.Range("A2").Select()
